$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-10 (columns B:G) with new computed values ---

# Row 2 (Q1)
$ws.Range("B2").Value = -0.02759765538578432
$ws.Range("C2").Value = 0.5993211969665078
$ws.Range("D2").Value = 0.7764809546060641
$ws.Range("E2").Value = 0.881181567332218
$ws.Range("F2").Value = 0.8895131873692134
$ws.Range("G2").Value = 51

# Row 3 (Q2)
$ws.Range("B3").Value = 0.09824497858667457
$ws.Range("C3").Value = 0.6346353876862429
$ws.Range("D3").Value = 0.9122639460283238
$ws.Range("E3").Value = 0.9551250944396361
$ws.Range("F3").Value = 0.9597043923121814
$ws.Range("G3").Value = 50

# Row 4 (Q3)
$ws.Range("B4").Value = -0.001319298517514361
$ws.Range("C4").Value = 0.6476906327130449
$ws.Range("D4").Value = 0.7146872890618677
$ws.Range("E4").Value = 0.8453917961879378
$ws.Range("F4").Value = 0.8541515268229536
$ws.Range("G4").Value = 49

# Row 5 (Q4)
$ws.Range("B5").Value = 0.1228341545534477
$ws.Range("C5").Value = 0.6897645129407867
$ws.Range("D5").Value = 0.839119603502874
$ws.Range("E5").Value = 0.9160347174113402
$ws.Range("F5").Value = 0.917367950390179
$ws.Range("G5").Value = 48

# Row 6 (Q5)
$ws.Range("B6").Value = 0.04626027039139741
$ws.Range("C6").Value = 0.6022393271950429
$ws.Range("D6").Value = 0.7064898076028564
$ws.Range("E6").Value = 0.8405294805078858
$ws.Range("F6").Value = 0.8483287965471219
$ws.Range("G6").Value = 47

# Row 7 (Q6)
$ws.Range("B7").Value = 0.1135060740419775
$ws.Range("C7").Value = 0.7098892706801739
$ws.Range("D7").Value = 0.8783339292678302
$ws.Range("E7").Value = 0.9371947125692879
$ws.Range("F7").Value = 0.940575637214629
$ws.Range("G7").Value = 46

# Row 8 (Q7)
$ws.Range("B8").Value = 0.06286716272461017
$ws.Range("C8").Value = 0.7192557759049942
$ws.Range("D8").Value = 0.9934777917745472
$ws.Range("E8").Value = 0.9967335610756504
$ws.Range("F8").Value = 1.005989427279836
$ws.Range("G8").Value = 45

# Row 9 (Q8)
$ws.Range("B9").Value = 0.1466813223655823
$ws.Range("C9").Value = 0.7501106339929442
$ws.Range("D9").Value = 0.9202103855659967
$ws.Range("E9").Value = 0.959275969450917
$ws.Range("F9").Value = 0.9589550867376859
$ws.Range("G9").Value = 44

# Row 10 (Q8 label kept) - note this row previously had no F10 value
$ws.Range("B10").Value = 0.07202425349984963
$ws.Range("C10").Value = 0.7064625585386706
$ws.Range("D10").Value = 0.795197501037768
$ws.Range("E10").Value = 0.8917384712110205
$ws.Range("F10").Value = 0.8993440776696643
$ws.Range("G10").Value = 43

# --- Add new row 11 (Q9), copying the formatting of row 10's label cell ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"

$ws.Range("B11").Value = 0.08705267901537297
$ws.Range("C11").Value = 0.729473168299468
$ws.Range("D11").Value = 0.8438328393602245
$ws.Range("E11").Value = 0.9186037444732219
$ws.Range("F11").Value = 0.9255544964036263
$ws.Range("G11").Value = 42
